$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

$ws.Range("D2").Value = '70.131.53'
$ws.Range("E2").Value = '  -1.07%  '

$ws.Range("D3").Value = '3.575.54'
$ws.Range("E3").Value = '  -1.64%  '

Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.22%  '

Set-TextValue $ws.Range("D5") '575.92'
$ws.Range("E5").Value = '  -2.76%  '

Set-TextValue $ws.Range("D6") '186.76'
$ws.Range("E6").Value = '  -3.85%  '

$ws.Range("D7").Value = '3.571.73'
$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("E8").Value = '  -3.37%  '

$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("E10").Value = '  +2.56%  '

$ws.Range("E11").Value = '  -3.00%  '

Set-TextValue $ws.Range("D12") '54.29'
$ws.Range("E12").Value = '  -5.91%  '

$ws.Range("E13").Value = '  -1.23%  '

Set-TextValue $ws.Range("D14") '9.58'
$ws.Range("E14").Value = '  -3.18%  '

$ws.Range("D15").Value = '4.149.23'
$ws.Range("E15").Value = '  -1.89%  '

Set-TextValue $ws.Range("D16") '19.68'
$ws.Range("E16").Value = '  -3.64%  '

$ws.Range("D17").Value = '3.583.01'
$ws.Range("E17").Value = '  -1.71%  '

$ws.Range("D18").Value = '70.034.73'
$ws.Range("E18").Value = '  -1.20%  '

$ws.Range("E19").Value = '  -2.12%  '

Set-TextValue $ws.Range("D20") '0.120'
$ws.Range("E20").Value = '  -1.29%  '

$ws.Range("E21").Value = '  -2.00%  '

Set-TextValue $ws.Range("D22") '495.86'
$ws.Range("E22").Value = '  +0.54%  '

Set-TextValue $ws.Range("D23") '19.43'
$ws.Range("E23").Value = '  +3.90%  '

Set-TextValue $ws.Range("D24") '5.05'
$ws.Range("E24").Value = '  -2.26%  '

Set-TextValue $ws.Range("D25") '96.56'
$ws.Range("E25").Value = '  +6.23%  '

Set-TextValue $ws.Range("D26") '4.37'
$ws.Range("E26").Value = '  -3.19%  '

Set-TextValue $ws.Range("D27") '11.58'
$ws.Range("E27").Value = '  +0.99%  '

Set-TextValue $ws.Range("D28") '2.99'
$ws.Range("E28").Value = '  -5.49%  '

Set-TextValue $ws.Range("D29") '9.35'
$ws.Range("E29").Value = '  -2.17%  '

Set-TextValue $ws.Range("D30") '7.73'
$ws.Range("E30").Value = '  -2.53%  '

Set-TextValue $ws.Range("D31") '31.67'
$ws.Range("E31").Value = '  -3.29%  '

Set-TextValue $ws.Range("D32") '12.84'
$ws.Range("E32").Value = '  +4.54%  '

Set-TextValue $ws.Range("D33") '65.69'
$ws.Range("E33").Value = '  -3.05%  '

$ws.Range("E34").Value = '  -4.69%  '

Set-TextValue $ws.Range("D35") '569.73'
$ws.Range("E35").Value = '  -7.16%  '

Set-TextValue $ws.Range("D36") '3.29'
$ws.Range("E36").Value = '  +11.76%  '

Set-TextValue $ws.Range("D37") '38.92'
$ws.Range("E37").Value = '  -3.42%  '

Set-TextValue $ws.Range("D38") '0.408'
$ws.Range("E38").Value = '  -0.46%  '

Set-TextValue $ws.Range("D39") '0.999'
$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("E40").Value = '  -5.95%  '

Set-TextValue $ws.Range("D41") '3.86'
$ws.Range("E41").Value = '  +13.80%  '

Set-TextValue $ws.Range("D42") '3.19'
$ws.Range("E42").Value = '  -1.19%  '

Set-TextValue $ws.Range("D43") '3.44'
$ws.Range("E43").Value = '  -3.68%  '

$ws.Range("E44").Value = '  -8.86%  '

Set-TextValue $ws.Range("D45") '3.05'
$ws.Range("E45").Value = '  -2.79%  '

Set-TextValue $ws.Range("D46") '0.0454'
$ws.Range("E46").Value = '  -0.99%  '

$ws.Range("D47").Value = '3.227.68'
$ws.Range("E47").Value = '  -2.90%  '

Set-TextValue $ws.Range("D48") '9.52'
$ws.Range("E48").Value = '  -1.60%  '

$ws.Range("E49").Value = '  -2.38%  '

Set-TextValue $ws.Range("D50") '1.50'
$ws.Range("E50").Value = '  +24.90%  '

Set-TextValue $ws.Range("D51") '0.998'
$ws.Range("E51").Value = '  -0.26%  '
